$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '27.154.80'
$ws.Range("E2").Value = '  +0.65%  '

$ws.Range("D3").Value = '1.906.52'
$ws.Range("E3").Value = '  +1.69%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.0000'
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.26'
$ws.Range("E5").Value = '  +0.59%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9990'
$ws.Range("E6").Value = '  -0.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5233'

$ws.Range("E8").Value = '  +3.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07243'
$ws.Range("E9").Value = '  +0.36%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.25'
$ws.Range("E10").Value = '  +2.19%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8962'
$ws.Range("E11").Value = '  +0.02%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07670'
$ws.Range("E12").Value = '  +1.95%  '

$ws.Range("D13").Value = '1.895.70'
$ws.Range("E13").Value = '  +0.80%  '

$ws.Range("E14").Value = '  -0.76%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.238'
$ws.Range("E15").Value = '  -0.22%  '

$ws.Range("E16").Value = '  -0.03%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008535'
$ws.Range("E17").Value = '  -0.06%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.54'
$ws.Range("E18").Value = '  +1.93%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9994'
$ws.Range("E19").Value = '  -0.11%  '

$ws.Range("D20").Value = '27.200.30'
$ws.Range("E20").Value = '  +0.70%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.080'
$ws.Range("E21").Value = '  +1.04%  '

$ws.Range("D22").Value = '2.137.65'
$ws.Range("E22").Value = '  +0.40%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.63'
$ws.Range("E23").Value = '  +2.20%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.436'
$ws.Range("E24").Value = '  +0.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.312'
$ws.Range("E25").Value = '  +10.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '145.68'
$ws.Range("E26").Value = '  -1.99%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.11'
$ws.Range("E27").Value = '  +1.07%  '

$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.728'
$ws.Range("E28").Value = '  -3.56%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.77'
$ws.Range("E29").Value = '  +1.13%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.964'
$ws.Range("E30").Value = '  +4.62%  '

$ws.Range("E31").Value = '  +1.65%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09206'
$ws.Range("E32").Value = '  +0.30%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05060'
$ws.Range("E33").Value = '  -0.95%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.244'
$ws.Range("E34").Value = '  +7.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7798'
$ws.Range("E35").Value = '  +3.89%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.983'
$ws.Range("E36").Value = '  +0.34%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.306'
$ws.Range("E37").Value = '  +2.14%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.607'
$ws.Range("E38").Value = '  +2.41%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5669'
$ws.Range("E39").Value = '  +0.84%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01994'
$ws.Range("E40").Value = '  -0.53%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.074'
$ws.Range("E41").Value = '  -0.43%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.023'
$ws.Range("E42").Value = '  +5.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.626'
$ws.Range("E43").Value = '  -0.30%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '118.67'
$ws.Range("E44").Value = '  +2.80%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1522'
$ws.Range("E45").Value = '  +2.83%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4853'
$ws.Range("E46").Value = '  +1.73%  '

$ws.Range("E47").Value = '  +1.32%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9985'
$ws.Range("E48").Value = '  -0.15%  '

$ws.Range("E49").Value = '  +1.87%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.57'
$ws.Range("E50").Value = '  +1.56%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '64.23'
$ws.Range("E51").Value = '  +1.57%  '

